# -------------------------------------------------------------------------
# Rebuild "Data.xlsx" from its original single-sheet layout into the new
# three-sheet layout described by the commit:
#   Sheet1 - small "id/msg" request-field table (new synthetic content)
#   Sheet2 - the original "sz00x / testing00x" rows, keeping only the
#            account (A), expected (C) and actual (D) columns
#   Sheet3 - same original rows but the password column (A) and a replaced
#            actual-result text (C) ("输入符合规则")
# -------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- capture the original Sheet1 data before we overwrite anything -------
$orig = $wb.Worksheets.Item(1)

$accounts  = @($orig.Range("A2").Value2, $orig.Range("A3").Value2, $orig.Range("A4").Value2, $orig.Range("A5").Value2)
$passwords = @($orig.Range("B2").Value2, $orig.Range("B3").Value2, $orig.Range("B4").Value2, $orig.Range("B5").Value2)

$headerAccount  = $orig.Range("A1").Value2
$headerPassword = $orig.Range("B1").Value2
$headerExpected = $orig.Range("C1").Value2
$headerActual   = $orig.Range("D1").Value2
$expectedMsg    = $orig.Range("C2").Value2

# ---- add the two extra worksheets (placed right after Sheet1) ------------
$sheet2 = $wb.Worksheets.Add([Type]::Missing, $orig, [Type]::Missing, [Type]::Missing)
$sheet3 = $wb.Worksheets.Add([Type]::Missing, $sheet2, [Type]::Missing, [Type]::Missing)

# =====================================================================
# Sheet1 - new small table of request field ids / messages
# =====================================================================
$orig.Cells.ClearContents()

$orig.Range("A1").Value = "id"
$orig.Range("B1").Value = "msg"
$orig.Range("A2").Value = "ctl00_holderLeft_txt_email"
$orig.Range("B2").Value = "tip_email"
$orig.Range("A3").Value = "ctl00_holderLeft_txt_userName"
$orig.Range("B3").Value = "tip_userName"

$orig.Columns.Item(1).ColumnWidth = 32.0357142857143
$orig.Columns.Item(2).ColumnWidth = 13.1607142857143

[void]$orig.Range("C12").Select()

# =====================================================================
# Sheet2 - account / expected / actual (columns A, C, D only)
# =====================================================================
$sheet2.Range("A1").Value = $headerAccount
$sheet2.Range("C1").Value = $headerExpected
$sheet2.Range("D1").Value = $headerActual

for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2
    $sheet2.Cells.Item($r, 1).Value = $accounts[$i]
    $sheet2.Cells.Item($r, 3).Value = $expectedMsg
}

$sheet2.Columns.Item(1).ColumnWidth = 5.78571428571429
$sheet2.Columns.Item(2).ColumnWidth = 12.2857142857143
$sheet2.Columns.Item(3).ColumnWidth = 26.1607142857143
$sheet2.Columns.Item(4).ColumnWidth = 12.2857142857143

[void]$sheet2.Range("D2:E5").Select()

# =====================================================================
# Sheet3 - password / expected / actual (columns A, C, D only)
# =====================================================================
$sheet3.Range("A1").Value = $headerPassword
$sheet3.Range("C1").Value = $headerExpected
$sheet3.Range("D1").Value = $headerActual

for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2
    $sheet3.Cells.Item($r, 1).Value = $passwords[$i]
    $sheet3.Cells.Item($r, 3).Value = "输入符合规则"
}

$sheet3.Columns.Item(1).ColumnWidth = 10.9107142857143
$sheet3.Columns.Item(3).ColumnWidth = 26.1607142857143
$sheet3.Columns.Item(4).ColumnWidth = 12.2857142857143

[void]$sheet3.Range("C21").Select()

# ---- Sheet3 ends up the active / selected tab -----------------------
[void]$sheet3.Activate()
